$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 669, shifting rows 669-681 down to 673-685
$ws.Range("A669:A672").EntireRow.Insert()

$ws.Cells.Item(669, 1).Value = 5
$ws.Cells.Item(669, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(669, 3).Value = "Maule"
$ws.Cells.Item(669, 4).Value = 44448
$ws.Cells.Item(669, 5).Value = 7
$ws.Cells.Item(669, 6).Value = 100112033
$ws.Cells.Item(669, 7).Value = "Lechuga"
$ws.Cells.Item(669, 8).Value = "Conconina(o)"
$ws.Cells.Item(669, 9).Value = "Segunda"
$ws.Cells.Item(669, 10).Value = 400
$ws.Cells.Item(669, 11).Value = 6000
$ws.Cells.Item(669, 12).Value = 6000
$ws.Cells.Item(669, 13).Value = 6000
$ws.Cells.Item(669, 14).Value = "$/caja 12 unidades"
$ws.Cells.Item(669, 15).Value = "Región del Maule"
$ws.Cells.Item(669, 16).Value = 500
$ws.Cells.Item(669, 17).Value = 12
$ws.Cells.Item(669, 18).Value = "Hortaliza"

$ws.Cells.Item(670, 1).Value = 5
$ws.Cells.Item(670, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(670, 3).Value = "Maule"
$ws.Cells.Item(670, 4).Value = 44448
$ws.Cells.Item(670, 5).Value = 7
$ws.Cells.Item(670, 6).Value = 100112033
$ws.Cells.Item(670, 7).Value = "Lechuga"
$ws.Cells.Item(670, 8).Value = "Escarola"
$ws.Cells.Item(670, 9).Value = "Primera"
$ws.Cells.Item(670, 10).Value = 600
$ws.Cells.Item(670, 11).Value = 9000
$ws.Cells.Item(670, 12).Value = 9000
$ws.Cells.Item(670, 13).Value = 9000
$ws.Cells.Item(670, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(670, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(670, 16).Value = 600
$ws.Cells.Item(670, 17).Value = 15
$ws.Cells.Item(670, 18).Value = "Hortaliza"

$ws.Cells.Item(671, 1).Value = 5
$ws.Cells.Item(671, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(671, 3).Value = "Maule"
$ws.Cells.Item(671, 4).Value = 44448
$ws.Cells.Item(671, 5).Value = 7
$ws.Cells.Item(671, 6).Value = 100112033
$ws.Cells.Item(671, 7).Value = "Lechuga"
$ws.Cells.Item(671, 8).Value = "Española"
$ws.Cells.Item(671, 9).Value = "Primera"
$ws.Cells.Item(671, 10).Value = 400
$ws.Cells.Item(671, 11).Value = 7000
$ws.Cells.Item(671, 12).Value = 7000
$ws.Cells.Item(671, 13).Value = 7000
$ws.Cells.Item(671, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(671, 15).Value = "Región del Maule"
$ws.Cells.Item(671, 16).Value = 389
$ws.Cells.Item(671, 17).Value = 18
$ws.Cells.Item(671, 18).Value = "Hortaliza"

$ws.Cells.Item(672, 1).Value = 5
$ws.Cells.Item(672, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(672, 3).Value = "Maule"
$ws.Cells.Item(672, 4).Value = 44448
$ws.Cells.Item(672, 5).Value = 7
$ws.Cells.Item(672, 6).Value = 100112033
$ws.Cells.Item(672, 7).Value = "Lechuga"
$ws.Cells.Item(672, 8).Value = "Marina"
$ws.Cells.Item(672, 9).Value = "Primera"
$ws.Cells.Item(672, 10).Value = 400
$ws.Cells.Item(672, 11).Value = 7000
$ws.Cells.Item(672, 12).Value = 7000
$ws.Cells.Item(672, 13).Value = 7000
$ws.Cells.Item(672, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(672, 15).Value = "Región del Maule"
$ws.Cells.Item(672, 16).Value = 389
$ws.Cells.Item(672, 17).Value = 18
$ws.Cells.Item(672, 18).Value = "Hortaliza"

